$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values for rows 2-21, columns A (x), B (y), C (z).
# The first 8 rows of the old "falling" sample were dropped, the
# remaining rows shifted up, and 8 new "struggle" rows were appended
# at the end so the overall extent (A1:C21) is unchanged.
$data = @(
    @(2, -0.6712930202484151, -3.392556905746461, 2.365111112594603),
    @(3, 0.5537151455879301, -4.361428594589236, 3.347476267814645),
    @(4, 2.452674245834349, -4.409869003295896, 5.084140586853025),
    @(5, 5.681596696376824, -3.525700151920317, 2.85166837722062),
    @(6, 1.737989616393984, -2.796510410308835, 1.393881118297584),
    @(7, -5.58917605876923, -7.849099040031435, 6.400659620761871),
    @(8, 2.166972637176496, -11.98566874265667, 9.183138275146462),
    @(9, -0.1177038192748863, -2.770210593938835, 3.860614097118379),
    @(10, 3.127950906753536, -4.12096252441406, 3.36216964721679),
    @(11, 2.318384975194932, -3.449181020259855, 1.645497059822083),
    @(12, 3.076034724712372, -2.935223340988159, 1.488467574119568),
    @(13, 4.307219874858854, -2.422872281074523, 1.113696080446245),
    @(14, 3.889325714111331, -2.331348705291748, 1.212168788909908),
    @(15, 4.565532350540161, -2.527925431728365, 0.07184810638427547),
    @(16, 4.525413262844082, -2.978561615943909, -0.2339930295944185),
    @(17, 3.761336714029313, -3.091410279273987, 0.3938730955123901),
    @(18, 4.025719761848449, -2.884652662277222, 0.07831997871398896),
    @(19, 3.922824805974961, -2.91312141418457, 0.03645527362823536),
    @(20, 4.07339057922363, -2.916565942764282, 0.1658158779144287),
    @(21, 3.406079268455506, -2.918866038322449, 0.1608816981315587)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
